$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row heights
$ws.Rows.Item(2).RowHeight = 27.75
$ws.Rows.Item(3).RowHeight = 33.75

# Font color fix: theme-black -> explicit black RGB for all cells using the "font 5" text style
$ws.Range("J2:J14").Font.Color = 0
$ws.Range("F4:F8").Font.Color = 0
$ws.Range("B9:C14").Font.Color = 0
$ws.Range("F9:G14").Font.Color = 0
$ws.Range("D9:D14").Font.Color = 0
$ws.Range("H9:H14").Font.Color = 0
$ws.Range("E12:E13").Font.Color = 0

# B3/C3 value + style change (match style used by B2/C2/F2/G2/F3/G3, i.e. s=11)
$ws.Range("B3").Value = 44656
$ws.Range("C3").Value = 44656
$ws.Range("B3:C3").NumberFormat = "d-mmm"
$ws.Range("B3:C3").HorizontalAlignment = -4152
$ws.Range("B3:C3").Font.Bold = $true
$ws.Range("B3:C3").Font.Size = 9
$ws.Range("B3:C3").Font.Color = 0
$ws.Range("B3:C3").Font.Name = "Calibri"

Write-Host "done"
